# Applies the "456a3b4" gh-pages data refresh to 苏州-漫展信息.xlsx
#
# Summary of the change (per the OOXML diff):
#   * On sheet "展览" and sheet "全部类型", a brand-new event row
#     ("苏州·OrangeOrange国潮&随机宅舞派对【免费活动】") is inserted right
#     after "苏州·OCG国潮动漫游戏嘉年华阿杰内场" and right before
#     "苏州·YoungComic动漫嘉年华". That shifts every following row down by
#     one and bumps the running index kept in column A.
#   * A handful of existing rows get their "想去人数" (F) / "最低票价" (G)
#     numbers refreshed to newer live values.

$wb = $excel.ActiveWorkbook

# Name of the event that marks the insertion point (new row goes directly
# above this one), and the full contents of the new row.
$anchorName = "苏州·YoungComic动漫嘉年华"

$newRowData = @{
    B = "2024-05-18"
    C = "苏州·OrangeOrange国潮&随机宅舞派对【免费活动】"
    D = "狮山路298号 金鹰国际购物中心(狮山路店)"
    E = "2024.05.18 13:00-05.18 17:00"
    F = 3
    G = 29
    H = "https://show.bilibili.com/platform/detail.html?id=83949"
    I = "//i1.hdslb.com/bfs/openplatform/202404/DOH6BK8i1712638105049.png"
}

# Refreshed numeric values, keyed by the event name (column C) so the same
# table drives both worksheets regardless of exact row numbers.
$numericUpdates = @{
    "苏州·X-party 国漫游戏嘉年华03"                                               = @{ F = 603 }
    "苏州·首届Redamancy动漫游戏嘉年华"                                            = @{ F = 1279 }
    "昆山·第十二届理想乡动漫游戏展"                                               = @{ F = 14222 }
    "苏州·I COME ACG动漫品牌博览会x中国国际动漫节cosplay超级盛典江苏赛区"          = @{ F = 16069; G = 29.9 }
    "苏州·萌动次元动漫游戏嘉年毕"                                                 = @{ F = 72 }
    "苏州·苏州湾动漫游戏嘉年华"                                                   = @{ F = 92 }
    "太仓·龙狮动漫嘉年华5.0"                                                      = @{ F = 32 }
    "常熟·CDW·动漫展03"                                                           = @{ F = 1235 }
    "苏州·AF动漫游戏嘉年华"                                                       = @{ F = 26 }
    "苏州·OCG国潮动漫游戏嘉年华"                                                  = @{ F = 6419 }
    "苏州·OCG国潮动漫游戏嘉年华阿杰内场"                                          = @{ F = 967 }
    "苏州·YoungComic动漫嘉年华"                                                   = @{ F = 1110 }
    "苏州·燃梦Project"                                                            = @{ F = 3 }
    "【会员购严选】苏州·Come in joy动漫国潮文化节"                                = @{ F = 5657 }
    "苏州·归离之缘原神only展"                                                     = @{ F = 89 }
    "苏州·白日梦想7.20全职猎人ONLY展"                                             = @{ F = 153 }
    "苏州·萤火国潮文化节动漫品牌博览会"                                           = @{ F = 4669 }
    "苏州·第十三届理想乡动漫展-同人创作者大会"                                    = @{ F = 9 }
}

$targetSheets = @("展览", "全部类型")

foreach ($sheetName in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    $lastRow = $ws.UsedRange.Rows.Count

    # Locate the anchor row (column C holds the event name).
    $insertRow = -1
    for ($r = 2; $r -le $lastRow; $r++) {
        if ($ws.Cells.Item($r, 3).Value() -eq $anchorName) {
            $insertRow = $r
            break
        }
    }

    # Insert a fresh row above the anchor; everything below shifts down.
    $ws.Rows.Item($insertRow).Insert()

    # Populate the new row's data cells (B..I); A is renumbered below.
    # B holds a plain "yyyy-mm-dd" label stored as text in this workbook
    # (every other row in the sheet is text too); a leading apostrophe
    # stops it from being auto-parsed into a date value/format.
    $ws.Cells.Item($insertRow, 2).Value = "'" + $newRowData.B
    $ws.Cells.Item($insertRow, 3).Value = $newRowData.C
    $ws.Cells.Item($insertRow, 4).Value = $newRowData.D
    $ws.Cells.Item($insertRow, 5).Value = $newRowData.E
    $ws.Cells.Item($insertRow, 6).Value = $newRowData.F
    $ws.Cells.Item($insertRow, 7).Value = $newRowData.G
    $ws.Cells.Item($insertRow, 8).Value = $newRowData.H
    $ws.Cells.Item($insertRow, 9).Value = $newRowData.I

    # Match the formatting used by every other "index" cell in column A
    # (bold, centered, thin box border) so the new cell reuses that style.
    $aCell = $ws.Cells.Item($insertRow, 1)
    $aCell.Font.Bold = $true
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160
    $aCell.Borders.LineStyle = 1

    $newLastRow = $lastRow + 1

    # Renumber column A (1-based running index, row 1 is the header = 0).
    for ($r = 2; $r -le $newLastRow; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }

    # Apply the refreshed F/G numbers, looked up by event name so the
    # shifted row positions don't matter.
    for ($r = 2; $r -le $newLastRow; $r++) {
        $name = $ws.Cells.Item($r, 3).Value()
        if ($numericUpdates.ContainsKey($name)) {
            $upd = $numericUpdates[$name]
            if ($upd.ContainsKey("F")) {
                $ws.Cells.Item($r, 6).Value = $upd.F
            }
            if ($upd.ContainsKey("G")) {
                $ws.Cells.Item($r, 7).Value = $upd.G
            }
        }
    }
}
